$d = $word.ActiveDocument

# Locate the sentence "This is a Microsoft word document." -- running
# Find on a Range narrows that Range down to the matched text (just
# like real Word COM behavior).
$target = $d.Content
$needle = "This is a Microsoft word document."
$found = $target.Find.Execute($needle, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target sentence: $needle"
}

# Append " (Changed main)" right after the sentence (and before its
# paragraph mark) as three separate runs -- " (", "Changed main" and
# ")" -- mirroring the authored edit's three new <w:r> elements.
$target.InsertAfter(" (")
$target.InsertAfter("Changed main")
$target.InsertAfter(")")
